$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: client location path changed to the new project folder
$ws.Range("B3").Value = "W:\Projects\בהת\176 יקותיאל אדם\קבצי עבודה\תחזיות_דמוגרפיות"

# B4: scenario name changed
$ws.Range("B4").Value = "with_project"

# B5: v_date is now a plain number instead of a text code
$ws.Range("B5").Value = 240929

# B6 keeps its value (output-by-version path) - unchanged content

# Update the active selection to B6
$ws.Range("B6").Select() | Out-Null
